$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Clear the "Pass" values that were prefilled in J2:J4 (DSL update for EB)
$ws.Range("J2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("J4").ClearContents()

# Move the active selection to J1, matching the new cursor position
$ws.Range("J1").Select()
